# Update graphs and tables: add a "Total" column (H) and a new "Global"
# summary row-pair (14-15) to the continent/sector N-studies table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "Total" header in H1 - copy the formatting (bold, centered,
#    thin border) from the neighbouring G1 header cell.
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "Total"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Column H ("Total") values for the existing continent rows 2-13.
# ---------------------------------------------------------------------
$totals = @{
    2  = "3674 (2400-5760)"
    3  = "1061 (605-1995)"
    4  = "21745 (14364-31884)"
    5  = "8868 (5002-15196)"
    6  = "5323 (3391-8104)"
    7  = "2251 (1220-4105)"
    8  = "13991 (9105-21466)"
    9  = "3762 (2089-7232)"
    10 = "20885 (14705-29783)"
    11 = "6764 (3999-11548)"
    12 = "5482 (3579-8202)"
    13 = "1922 (1047-3441)"
}
foreach ($r in $totals.Keys) {
    $ws.Range("H$r").Value = $totals[$r]
}

# ---------------------------------------------------------------------
# 3. New "Global" summary rows: row 14 (D&A) and row 15 (Other),
#    mirroring the layout of the per-continent row pairs above.
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Global"
$ws.Range("B14").Value = "D&A"
$ws.Range("C14").Value = "5529 (2957-9755)"
$ws.Range("D14").Value = "7729 (4487-14423)"
$ws.Range("E14").Value = "3783 (2201-7050)"
$ws.Range("F14").Value = "6623 (3777-12384)"
$ws.Range("G14").Value = "19622 (10796-34361)"
$ws.Range("H14").Value = "57366 (38371-85227)"

$ws.Range("B15").Value = "Other"
$ws.Range("C15").Value = "7078 (2838-15628)"
$ws.Range("D15").Value = "4457 (2036-11717)"
$ws.Range("E15").Value = "2524 (1325-5178)"
$ws.Range("F15").Value = "6182 (3278-12777)"
$ws.Range("G15").Value = "15366 (7724-31305)"
$ws.Range("H15").Value = "44373 (25594-78626)"

# Merge the A14:A15 continent-name cell exactly like A2:A3, A4:A5, etc.
# (Do this before copying the label formatting below - merging after
# formatting disturbs the per-cell border resolution.)
$ws.Range("A14:A15").Merge()

# Copy the bold / bordered / centered "label" formatting from the
# Oceania row pair (A12:B13, the last existing continent rows) onto
# the new Global row pair.
$ws.Range("A12").Copy()
$ws.Range("A14").PasteSpecial(-4122)

$ws.Range("A13").Copy()
$ws.Range("A15").PasteSpecial(-4122)

$ws.Range("B12").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial(-4122)
